$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.121747374534607
$ws.Range("B1").Value = 2.267987012863159
$ws.Range("C1").Value = 10.37368297576904
$ws.Range("D1").Value = 1.848957538604736
$ws.Range("E1").Value = 1.28792417049408
